# Team Everyday Attendance - add a new day's row (07-08-2023) below the
# last existing entry, mirroring the previous row's PRESENT/ABSENT pattern,
# plus "No Response" comments on the ABSENT cells for the last three people.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The prior row (row 4, 04-Aug-2023) already carries the exact style/format
# (date number format) and PRESENT/ABSENT pattern that the new row needs, so
# copy it down to row 5 and then just correct the date.
$ws.Range("A4:I4").Copy($ws.Range("A5:I5"))
$ws.Range("A5").Value = 45145

# Renuka's "No Response" notes on the absentees for the new day.
$note = "RENUKA:`nNo Response`n"
$ws.Range("G5").AddComment($note)
$ws.Range("H5").AddComment($note)
$ws.Range("I5").AddComment($note)

# Leave the selection where Excel would land after entering the new row.
$ws.Range("I5").Select() | Out-Null
